# Apply cryptos list update (prices/volumes refresh + two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.600.66"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "1.895.89"

$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "'327.02"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").Value = "'0.4597"
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("E8").Value = "  -1.68%  "

$ws.Range("D9").Value = "'46.83"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").Value = "'0.07878"
$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  +2.60%  "

$ws.Range("D12").Value = "'21.71"
$ws.Range("E12").Value = "  -2.97%  "

$ws.Range("D13").Value = "1.926.22"
$ws.Range("E13").Value = "  +3.40%  "

$ws.Range("D14").Value = "'7.099"
$ws.Range("E14").Value = "  +2.01%  "

$ws.Range("D15").Value = "'5.717"
$ws.Range("E15").Value = "  -0.70%  "

$ws.Range("D16").Value = "'0.06966"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "'87.43"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "'1.007"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").Value = "'0.00001005"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").Value = "'17.26"
$ws.Range("E20").Value = "  +1.83%  "

$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").Value = "28.628.43"
$ws.Range("E22").Value = "  -0.77%  "

$ws.Range("D23").Value = "'5.322"
$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("D24").Value = "'11.03"
$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("D25").Value = "2.129.77"
$ws.Range("E25").Value = "  +2.08%  "

$ws.Range("D26").Value = "'2.060"
$ws.Range("E26").Value = "  -2.63%  "

$ws.Range("D27").Value = "'154.86"
$ws.Range("E27").Value = "  +0.69%  "

$ws.Range("D28").Value = "'19.34"
$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("D29").Value = "'5.857"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'1.941"
$ws.Range("E30").Value = "  -3.41%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'118.59"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").Value = "'0.09329"
$ws.Range("E32").Value = "  -0.59%  "

$ws.Range("D33").Value = "'0.9300"
$ws.Range("E33").Value = "  -1.30%  "

$ws.Range("D34").Value = "'5.306"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").Value = "'1.335"
$ws.Range("E35").Value = "  -1.34%  "

$ws.Range("D36").Value = "'3.266"
$ws.Range("E36").Value = "  -2.72%  "

$ws.Range("D37").Value = "'0.05750"
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").Value = "'1.156"
$ws.Range("E38").Value = "  +0.64%  "

$ws.Range("D39").Value = "'0.02072"
$ws.Range("E39").Value = "  -2.39%  "

$ws.Range("D40").Value = "'7.758"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("D41").Value = "'0.5639"
$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("D42").Value = "'0.1786"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("E43").Value = "  -2.18%  "

$ws.Range("D44").Value = "'2.218"
$ws.Range("E44").Value = "  +4.75%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.07170"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'11.75"
$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("D47").Value = "'0.5340"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'1.116"
$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("D49").Value = "'1.835"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").Value = "'113.13"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("D51").Value = "'2.458"
$ws.Range("E51").Value = "  +3.96%  "

